$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 10.752749
$ws.Range("H2").Value = 32.258247
$ws.Range("I2").Value = 0.2551491597938751
$ws.Range("J2").Value = 0.2551491597938751
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 93.85711466666667
$ws.Range("N2").Value = 281.571344
$ws.Range("Q2").Value = 1009.221995874885
$ws.Range("R2").Value = 9082.997962873967
$ws.Range("S2").Value = 0.2551491597938751
$ws.Range("T2").Value = 0.2551491597938751

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 19.86979433333333
$ws.Range("H3").Value = 59.60938299999999
$ws.Range("I3").Value = 0.4714851364453034
$ws.Range("J3").Value = 0.4714851364453035
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 93.85711466666667
$ws.Range("N3").Value = 281.571344
$ws.Range("Q3").Value = 1864.92156514675
$ws.Range("R3").Value = 16784.29408632075
$ws.Range("S3").Value = 0.4714851364453034
$ws.Range("T3").Value = 0.4714851364453035

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 11.52044866666667
$ws.Range("H4").Value = 34.561346
$ws.Range("I4").Value = 0.2733657037608214
$ws.Range("J4").Value = 0.2733657037608214
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 93.85711466666667
$ws.Range("N4").Value = 281.571344
$ws.Range("Q4").Value = 1081.276071518781
$ws.Range("R4").Value = 9731.484643669024
$ws.Range("S4").Value = 0.2733657037608214
$ws.Range("T4").Value = 0.2733657037608214
